# Auto commit at 2025-09-10 7:40:33.70
#
# "Metrics" sheet: refresh the monthly/yearly/total figures in column B
# (rows 2-13) with the latest numbers, and leave the selection where the
# user finished editing (B2:B13, active cell B2).
#
# "today" sheet: it pulls these same figures in via formulas
# (=Metrics!B2, =Metrics!B3, ...) plus derived E/F columns, so it
# recalculates automatically once Metrics changes. We only need to move
# its selection to the cell the user ended up on (E7); it remains the
# active sheet/tab, same as before the edit.

$wb = $excel.ActiveWorkbook

$wsMetrics = $wb.Worksheets.Item("Metrics")

$wsMetrics.Range("B2").Value  = 141078.84
$wsMetrics.Range("B3").Value  = 114479.68000000002
$wsMetrics.Range("B4").Value  = 44695.100000000006
$wsMetrics.Range("B5").Value  = 5555
$wsMetrics.Range("B6").Value  = 4060329.7199999997
$wsMetrics.Range("B7").Value  = 3442007.1599999992
$wsMetrics.Range("B8").Value  = 1174060.78
$wsMetrics.Range("B9").Value  = 156715
$wsMetrics.Range("B10").Value = 32525653.520999826
$wsMetrics.Range("B11").Value = 19471877.230000004
$wsMetrics.Range("B12").Value = 11455769.670000002
$wsMetrics.Range("B13").Value = 1254342

# Reflect the new selection on the Metrics sheet.
$wsMetrics.Activate() | Out-Null
$wsMetrics.Range("B2:B13").Select() | Out-Null

# Restore focus/selection on the "today" sheet (it was, and remains, the
# active tab).
$wsToday = $wb.Worksheets.Item("today")
$wsToday.Activate() | Out-Null
$wsToday.Range("E7").Select() | Out-Null
